$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 314, shifting the
# existing rows 314-341 down to 316-343 (matching the dimension growing
# from A1:T341 to A1:T343).
$ws.Rows("314:315").Insert()

# Populate the two newly inserted rows with the new weekly price records.
$ws.Cells.Item(314, 1).Value = 4
$ws.Cells.Item(314, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(314, 3).Value = "Los Lagos"
$ws.Cells.Item(314, 4).Value = 44826
$ws.Cells.Item(314, 5).Value = 10
$ws.Cells.Item(314, 6).Value = "Fruta"
$ws.Cells.Item(314, 7).Value = 100104
$ws.Cells.Item(314, 8).Value = "Frutos de pepita"
$ws.Cells.Item(314, 9).Value = 100104005
$ws.Cells.Item(314, 10).Value = "Pera"
$ws.Cells.Item(314, 11).Value = "Packham's Triumph"
$ws.Cells.Item(314, 12).Value = "Primera"
$ws.Cells.Item(314, 13).Value = 400
$ws.Cells.Item(314, 14).Value = 15000
$ws.Cells.Item(314, 15).Value = 16000
$ws.Cells.Item(314, 16).Value = 15500
$ws.Cells.Item(314, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(314, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(314, 19).Value = 1033
$ws.Cells.Item(314, 20).Value = 15

$ws.Cells.Item(315, 1).Value = 4
$ws.Cells.Item(315, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(315, 3).Value = "Los Lagos"
$ws.Cells.Item(315, 4).Value = 44826
$ws.Cells.Item(315, 5).Value = 10
$ws.Cells.Item(315, 6).Value = "Fruta"
$ws.Cells.Item(315, 7).Value = 100104
$ws.Cells.Item(315, 8).Value = "Frutos de pepita"
$ws.Cells.Item(315, 9).Value = 100104005
$ws.Cells.Item(315, 10).Value = "Pera"
$ws.Cells.Item(315, 11).Value = "Packham's Triumph"
$ws.Cells.Item(315, 12).Value = "Segunda"
$ws.Cells.Item(315, 13).Value = 200
$ws.Cells.Item(315, 14).Value = 13000
$ws.Cells.Item(315, 15).Value = 13000
$ws.Cells.Item(315, 16).Value = 13000
$ws.Cells.Item(315, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(315, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(315, 19).Value = 867
$ws.Cells.Item(315, 20).Value = 15
